# v1.2 close user home review
$wb = $excel.ActiveWorkbook

$wsReviews = $wb.Worksheets.Item("LH-TC-USERHOME-Reviews")
$wsHistory = $wb.Worksheets.Item("Version History")

# Close out the review status on row 5 (I5 / J5): Open -> Closed
$wsReviews.Range("I5").Value = "Closed"
$wsReviews.Range("J5").Value = "Closed"
$wsReviews.Range("J5").Select()

# Add a new Version History row documenting the v1.2 close-out.
# Copy formatting from the row above (row 3) so the new row matches the
# existing table styling, then fill in the values.
$wsHistory.Range("A3:D3").Copy()
$wsHistory.Range("A4:D4").PasteSpecial(-4122)
$wsHistory.Rows.Item(4).RowHeight = 30
$wsHistory.Range("A4").Value = "v1.2"
$wsHistory.Range("B4").Value = "Ahmed Abuzaid"
$wsHistory.Range("C4").Value = "close user home review "
$wsHistory.Range("D4").Value = "14/5/2025"

# Grow the "Table1" ListObject so the new row is included in the table.
$lo = $wsHistory.ListObjects.Item(1)
$lo.Resize($wsHistory.Range("A1:D4"))

$wsHistory.Range("C4").Select()
